$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit rotates the observation-record data among rows 9, 10, 12, 15 and 16
# (row 11, 13, 14 and all other rows are untouched), and updates a few cells
# that gain/lose values as a result (blank placeholder columns such as
# J/K/L/M/N/AC/AF move with the record they belong to).

# ---- Row 9 (becomes the old row 10 record: Kortskaftad argspik) ----
$ws.Range("A9").Value = 111541121
$ws.Range("B9").Value = 79444
$ws.Range("E9").Value = 1049
$ws.Range("F9").Value = "Kortskaftad ärgspik"
$ws.Range("G9").Value = "Microcalicium ahlneri"
$ws.Range("H9").Value = "Tibell"
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = ""
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = ""
$ws.Range("Q9").Value = 693460.9606228607
$ws.Range("R9").Value = 6551521.405726598
$ws.Range("AC9").Value = ""
$ws.Range("AF9").Value = ""
$ws.Range("AO9").Value = "silverstubbe av tall"

# ---- Row 10 (becomes the old row 16 record: Bronshjon / farska gnagspar) ----
$ws.Range("A10").Value = 111541128
$ws.Range("B10").Value = 5113
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 100526
$ws.Range("F10").Value = "Bronshjon"
$ws.Range("G10").Value = "Callidium coriaceum"
$ws.Range("H10").Value = "Paykull, 1800"
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = "färska gnagspår"
$ws.Range("N10").Value = ""
$ws.Range("Q10").Value = 693570.8046739453
$ws.Range("R10").Value = 6551451.742365629
$ws.Range("AF10").Value = ""
$ws.Range("AO10").Value = "torrgran"

# ---- Row 12 (becomes the old row 15 record: Stor revmossa) ----
$ws.Range("A12").Value = 111541118
$ws.Range("B12").Value = 94851
$ws.Range("E12").Value = 2569
$ws.Range("F12").Value = "Stor revmossa"
$ws.Range("G12").Value = "Bazzania trilobata"
$ws.Range("H12").Value = "(L.) Gray"
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("Q12").Value = 693461.6376634488
$ws.Range("R12").Value = 6551559.049034445
$ws.Range("AF12").Value = ""
$ws.Range("AO12").Value = ""

# ---- Row 15 (becomes the old row 9 record: Reliktbock / gammeltall) ----
$ws.Range("A15").Value = 111541119
$ws.Range("B15").Value = 5426
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 101410
$ws.Range("F15").Value = "Reliktbock"
$ws.Range("G15").Value = "Nothorhina muricata"
$ws.Range("H15").Value = "(Dalman, 1817)"
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = "äldre gnagspår"
$ws.Range("N15").Value = ""
$ws.Range("Q15").Value = 693467.6220677271
$ws.Range("R15").Value = 6551532.561666255
$ws.Range("AC15").Value = "En gammal tall med kläckhål här och var. Om det är färskt eller gammalt är svårt sia om."
$ws.Range("AF15").Value = ""
$ws.Range("AO15").Value = "gammeltall"

# ---- Row 16 (becomes the old row 12 record: Bronshjon / aldre gnagspar) ----
$ws.Range("A16").Value = 111541129
$ws.Range("M16").Value = "äldre gnagspår"
$ws.Range("Q16").Value = 693328.6441019299
$ws.Range("R16").Value = 6551545.628735202
